$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2026-01-02 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-03 Saturday", 2) | Out-Null

# Update each multiplication fact cell-by-cell (table has repeated values, so
# scope each Find/Replace to its own cell Range with wdReplaceOne).
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$cell.Range.Find.Execute("543×8=4344", $true, $false, $false, $false, $false, $true, 1, $false, "148×4=592", 1) | Out-Null

$cell = $tbl.Cell(1, 2)
$cell.Range.Find.Execute("670×2=1340", $true, $false, $false, $false, $false, $true, 1, $false, "361×5=1805", 1) | Out-Null

$cell = $tbl.Cell(1, 3)
$cell.Range.Find.Execute("651×3=1953", $true, $false, $false, $false, $false, $true, 1, $false, "738×8=5904", 1) | Out-Null

$cell = $tbl.Cell(1, 4)
$cell.Range.Find.Execute("444×5=2220", $true, $false, $false, $false, $false, $true, 1, $false, "183×6=1098", 1) | Out-Null

$cell = $tbl.Cell(1, 5)
$cell.Range.Find.Execute("454×6=2724", $true, $false, $false, $false, $false, $true, 1, $false, "970×3=2910", 1) | Out-Null

$cell = $tbl.Cell(5, 1)
$cell.Range.Find.Execute("673×3=2019", $true, $false, $false, $false, $false, $true, 1, $false, "497×5=2485", 1) | Out-Null

$cell = $tbl.Cell(5, 2)
$cell.Range.Find.Execute("494×7=3458", $true, $false, $false, $false, $false, $true, 1, $false, "297×9=2673", 1) | Out-Null

$cell = $tbl.Cell(5, 3)
$cell.Range.Find.Execute("670×2=1340", $true, $false, $false, $false, $false, $true, 1, $false, "828×3=2484", 1) | Out-Null

$cell = $tbl.Cell(5, 4)
$cell.Range.Find.Execute("719×5=3595", $true, $false, $false, $false, $false, $true, 1, $false, "147×7=1029", 1) | Out-Null

$cell = $tbl.Cell(5, 5)
$cell.Range.Find.Execute("157×3=471", $true, $false, $false, $false, $false, $true, 1, $false, "279×9=2511", 1) | Out-Null

$cell = $tbl.Cell(10, 1)
$cell.Range.Find.Execute("930×8=7440", $true, $false, $false, $false, $false, $true, 1, $false, "387×9=3483", 1) | Out-Null

$cell = $tbl.Cell(10, 2)
$cell.Range.Find.Execute("602×7=4214", $true, $false, $false, $false, $false, $true, 1, $false, "244×8=1952", 1) | Out-Null

$cell = $tbl.Cell(10, 3)
$cell.Range.Find.Execute("696×3=2088", $true, $false, $false, $false, $false, $true, 1, $false, "466×8=3728", 1) | Out-Null

$cell = $tbl.Cell(10, 4)
$cell.Range.Find.Execute("913×8=7304", $true, $false, $false, $false, $false, $true, 1, $false, "640×6=3840", 1) | Out-Null

$cell = $tbl.Cell(10, 5)
$cell.Range.Find.Execute("755×9=6795", $true, $false, $false, $false, $false, $true, 1, $false, "392×5=1960", 1) | Out-Null

$cell = $tbl.Cell(15, 1)
$cell.Range.Find.Execute("916×7=6412", $true, $false, $false, $false, $false, $true, 1, $false, "755×3=2265", 1) | Out-Null

$cell = $tbl.Cell(15, 2)
$cell.Range.Find.Execute("849×4=3396", $true, $false, $false, $false, $false, $true, 1, $false, "880×2=1760", 1) | Out-Null

$cell = $tbl.Cell(15, 3)
$cell.Range.Find.Execute("567×4=2268", $true, $false, $false, $false, $false, $true, 1, $false, "182×7=1274", 1) | Out-Null

$cell = $tbl.Cell(15, 4)
$cell.Range.Find.Execute("547×8=4376", $true, $false, $false, $false, $false, $true, 1, $false, "318×6=1908", 1) | Out-Null

$cell = $tbl.Cell(15, 5)
$cell.Range.Find.Execute("120×5=600", $true, $false, $false, $false, $false, $true, 1, $false, "294×6=1764", 1) | Out-Null

$cell = $tbl.Cell(20, 1)
$cell.Range.Find.Execute("365×9=3285", $true, $false, $false, $false, $false, $true, 1, $false, "379×3=1137", 1) | Out-Null

$cell = $tbl.Cell(20, 2)
$cell.Range.Find.Execute("335×4=1340", $true, $false, $false, $false, $false, $true, 1, $false, "770×6=4620", 1) | Out-Null

$cell = $tbl.Cell(20, 3)
$cell.Range.Find.Execute("993×3=2979", $true, $false, $false, $false, $false, $true, 1, $false, "631×2=1262", 1) | Out-Null

$cell = $tbl.Cell(20, 4)
$cell.Range.Find.Execute("451×5=2255", $true, $false, $false, $false, $false, $true, 1, $false, "670×8=5360", 1) | Out-Null

$cell = $tbl.Cell(20, 5)
$cell.Range.Find.Execute("354×2=708", $true, $false, $false, $false, $false, $true, 1, $false, "561×7=3927", 1) | Out-Null
